$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: was inline string "100" -> now numeric 22
$ws.Range("E2").Value = 22

# E5: was numeric 50 -> now numeric 48
$ws.Range("E5").Value = 48
